# Update gh-pages output values as generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) — first block of changes
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F12").Value = 578
$wsExpo.Range("F13").Value = 268
$wsExpo.Range("F18").Value = 9529
$wsExpo.Range("F22").Value = 12047
$wsExpo.Range("F28").Value = 2712
$wsExpo.Range("F32").Value = 1012
$wsExpo.Range("F37").Value = 1106

# Sheet "全部类型" (All types) — second block of changes
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F17").Value = 578
$wsAll.Range("F18").Value = 268
$wsAll.Range("F22").Value = 9529
$wsAll.Range("F26").Value = 12047
$wsAll.Range("F34").Value = 2712
